$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.445.99"
$ws.Range("E2").Value = "  -0.24%  "

$ws.Range("D3").Value = "3.524.79"
$ws.Range("E3").Value = "  -0.34%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "613.01"
$ws.Range("E5").Value = "  -0.02%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "151.38"
$ws.Range("E6").Value = "  -1.84%  "

$ws.Range("D7").Value = "3.523.57"
$ws.Range("E7").Value = "  -0.28%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.482"
$ws.Range("E9").Value = "  -0.78%  "

$ws.Range("E10").Value = "  -1.04%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.12"
$ws.Range("E11").Value = "  +3.27%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.425"
$ws.Range("E12").Value = "  -1.45%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000220"
$ws.Range("E13").Value = "  -1.21%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.08"
$ws.Range("E14").Value = "  -0.47%  "

$ws.Range("D15").Value = "4.120.33"
$ws.Range("E15").Value = "  -0.27%  "

$ws.Range("D16").Value = "3.531.34"
$ws.Range("E16").Value = "  +0.01%  "

$ws.Range("D17").Value = "67.407.55"
$ws.Range("E17").Value = "  -0.23%  "

$ws.Range("E18").Value = "  -0.02%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.39"
$ws.Range("E19").Value = "  +0.06%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.28"
$ws.Range("E20").Value = "  -2.03%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "444.88"
$ws.Range("E21").Value = "  -2.11%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.54"
$ws.Range("E22").Value = "  +1.03%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.624"
$ws.Range("E23").Value = "  -3.03%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "77.49"
$ws.Range("E24").Value = "  -1.56%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000131"
$ws.Range("E25").Value = "  +9.19%  "

$ws.Range("D26").Value = "3.665.46"
$ws.Range("E26").Value = "  -0.17%  "

$ws.Range("E27").Value = "  +0.01%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.24"
$ws.Range("E28").Value = "  -2.51%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.46"
$ws.Range("E29").Value = "  +1.03%  "

$ws.Range("E30").Value = "  -2.32%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.55"
$ws.Range("E31").Value = "  -8.74%  "

$ws.Range("E32").Value = "  +0.01%  "

$ws.Range("E33").Value = "  +3.99%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.85"
$ws.Range("E34").Value = "  -0.73%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.17"
$ws.Range("E35").Value = "  -1.05%  "

$ws.Range("D36").Value = "3.516.31"
$ws.Range("E36").Value = "  -0.40%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.84"
$ws.Range("E37").Value = "  -3.93%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.00"
$ws.Range("E38").Value = "  -0.29%  "

$ws.Range("E39").Value = "  -0.02%  "

$ws.Range("E40").Value = "  +0.04%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "177.14"
$ws.Range("E41").Value = "  +3.18%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.15"
$ws.Range("E42").Value = "  +1.61%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0883"
$ws.Range("E43").Value = "  +0.17%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.42"
$ws.Range("E44").Value = "  -3.65%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.880"
$ws.Range("E45").Value = "  -1.39%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "28.18"
$ws.Range("E46").Value = "  -4.64%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "45.14"
$ws.Range("E47").Value = "  -1.38%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.63"
$ws.Range("E48").Value = "  -0.66%  "

$ws.Range("E49").Value = "  +0.91%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.59"
$ws.Range("E50").Value = "  -1.03%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.995"
$ws.Range("E51").Value = "  -2.84%  "

